# Auto-generated script applying the crypto price/volume update diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.439.88'
$ws.Range("E2").Value = '  +1.32%  '
$ws.Range("D3").Value = '2.275.96'
$ws.Range("E3").Value = '  +2.42%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '250.89'
$ws.Range("E5").Value = '  -0.07%  '
$ws.Range("E6").Value = '  +1.38%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '72.16'
$ws.Range("E7").Value = '  +5.86%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.641'
$ws.Range("E9").Value = '  +1.02%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.97'
$ws.Range("E10").Value = '  -1.85%  '
$ws.Range("E11").Value = '  +2.64%  '
$ws.Range("E12").Value = '  -1.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.26'
$ws.Range("E13").Value = '  +2.30%  '
$ws.Range("E14").Value = '  +1.89%  '
$ws.Range("D15").Value = '2.617.73'
$ws.Range("E15").Value = '  +2.68%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '15.01'
$ws.Range("E16").Value = '  +2.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.877'
$ws.Range("E17").Value = '  +0.12%  '
$ws.Range("D18").Value = '2.275.93'
$ws.Range("E18").Value = '  +3.09%  '
$ws.Range("D19").Value = '42.395.49'
$ws.Range("E19").Value = '  +1.41%  '
$ws.Range("D20").Value = '0.0₃0996'
$ws.Range("E20").Value = '  +3.38%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.29'
$ws.Range("E21").Value = '  +0.77%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '72.07'
$ws.Range("E22").Value = '  -1.02%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '235.03'
$ws.Range("E23").Value = '  +1.35%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.24'
$ws.Range("E24").Value = '  +7.93%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.91'
$ws.Range("E25").Value = '  +0.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.49'
$ws.Range("E26").Value = '  +0.93%  '
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.43'
$ws.Range("E28").Value = '  +0.46%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '3.66'
$ws.Range("E29").Value = '  -1.05%  '
$ws.Range("E30").Value = '  +2.33%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '167.21'
$ws.Range("E31").Value = '  -0.08%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.06'
$ws.Range("E32").Value = '  +2.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.41'
$ws.Range("E33").Value = '  +7.48%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.126'
$ws.Range("E34").Value = '  +4.66%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0813'
$ws.Range("E35").Value = '  +1.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.13'
$ws.Range("E36").Value = '  +21.92%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.126'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.74'
$ws.Range("E38").Value = '  +14.89%  '
$ws.Range("E39").Value = '  +1.85%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0306'
$ws.Range("E40").Value = '  -0.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '13.92'
$ws.Range("E41").Value = '  +14.61%  '
$ws.Range("E42").Value = '  +3.31%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.91'
$ws.Range("E43").Value = '  +4.74%  '
$ws.Range("E44").Value = '  +6.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.15'
$ws.Range("E45").Value = '  +6.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '61.47'
$ws.Range("E46").Value = '  -1.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.90'
$ws.Range("E47").Value = '  -4.77%  '
$ws.Range("E48").Value = '  +3.07%  '
$ws.Range("E49").Value = '  +0.39%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.18'
$ws.Range("E50").Value = '  +0.92%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '97.38'
$ws.Range("E51").Value = '  +4.64%  '
